# Reorder the sheet tabs so "review_info" comes before "hotel_info"
$wb = $excel.ActiveWorkbook
$reviewInfo = $wb.Worksheets.Item("review_info")
$hotelInfo  = $wb.Worksheets.Item("hotel_info")
$reviewInfo.Move($hotelInfo)

# Re-fetch the hotel_info worksheet by name now that the tab order changed
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Add a new "State" column to hotel_info, right after "Hotel_Name"
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"
